$wb = $excel.ActiveWorkbook

# --- Sheet: optimization_parameters ---
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the duplicated "value" header cells in C1:F1 (row span becomes 1:5)
$ws.Range("C1:F1").ClearContents()

# Insert a new row after row 8 (TolX) to make room for the new
# "production_function" / "L_curve" parameter rows.
$ws.Rows.Item(9).Insert()

# Old row 8 ("Model" / "Sigmoid") becomes "production_function" / "Sigmoid"
$ws.Range("A8").Value2 = "production_function"

# New row 9: "L_curve" = 1 (same number format as the other small numeric
# parameters above, e.g. alpha in B2)
$ws.Range("A9").Value2 = "L_curve"
$ws.Range("B9").Value2 = 1
$ws.Range("B9").NumberFormat = "0.00E+00"

# The old "Deletion" row (now shifted down to row 17 after the insert above)
# is removed entirely.
$ws.Rows.Item(17).Delete()

# --- Active sheet / window changes ---
# Previously "optimization_diagnostics" (index 14) was the active/selected
# tab; now it's "optimization_parameters" (index 7), with a specific
# selection.
$ws.Activate()
$ws.Range("C1:H7").Select()
